# Adding Date of Birth to participant signups, part 3
#
# Inserts a new "DOB" column ahead of the existing "Age" column (old column G),
# shifting every column from G onward one position to the right, then selects
# cell G2 (the new column's data cell for the single participant row) to match
# the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at G; everything from G onward (Age, Gender, ...,
# CampPrefs) shifts right by one column (G->H, H->I, ..., AP->AQ).
$ws.Columns("G").Insert()

# Give the newly inserted column its header - this also creates the new
# shared-string entry "DOB" right where the old "Age" string used to be.
$ws.Range("G1").Value = "DOB"

# Match the author's final selection.
$ws.Range("G2").Select() | Out-Null
